$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cluster labels in data rows 2-10 with new combination grid (ECs, FAPs, sCs)
# and updated numeric stats per Dr Hou advice.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Inhba"
$ws.Range("C2").Value = "Acvr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.123204333333334
$ws.Range("H2").Value = 12.369613
$ws.Range("I2").Value = 0.2909967288544799
$ws.Range("J2").Value = 0.2909967288544799
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.057757666666666
$ws.Range("N2").Value = 15.173273
$ws.Range("O2").Value = 0.173378811020062
$ws.Range("P2").Value = 0.173378811020062
$ws.Range("Q2").Value = 20.85416832814989
$ws.Range("R2").Value = 187.687514953349
$ws.Range("S2").Value = 0.05045266685951709
$ws.Range("T2").Value = 0.0504526668595171

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Inhba"
$ws.Range("C3").Value = "Acvr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.123204333333334
$ws.Range("H3").Value = 12.369613
$ws.Range("I3").Value = 0.2909967288544799
$ws.Range("J3").Value = 0.2909967288544799
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.247411
$ws.Range("N3").Value = 42.742233
$ws.Range("O3").Value = 0.4883980890531961
$ws.Range("P3").Value = 0.4883980890531961
$ws.Range("Q3").Value = 58.744986773981
$ws.Range("R3").Value = 528.704880965829
$ws.Range("S3").Value = 0.1421222462932591
$ws.Range("T3").Value = 0.1421222462932591

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Inhba"
$ws.Range("C4").Value = "Acvr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 4.123204333333334
$ws.Range("H4").Value = 12.369613
$ws.Range("I4").Value = 0.2909967288544799
$ws.Range("J4").Value = 0.2909967288544799
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.866548666666667
$ws.Range("N4").Value = 29.599646
$ws.Range("O4").Value = 0.3382230999267418
$ws.Range("P4").Value = 0.3382230999267418
$ws.Range("Q4").Value = 40.68179621744422
$ws.Range("R4").Value = 366.136165956998
$ws.Range("S4").Value = 0.09842181570170376
$ws.Range("T4").Value = 0.09842181570170376

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Inhba"
$ws.Range("C5").Value = "Acvr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.433639666666666
$ws.Range("H5").Value = 25.300919
$ws.Range("I5").Value = 0.5952073574179045
$ws.Range("J5").Value = 0.5952073574179045
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.057757666666666
$ws.Range("N5").Value = 15.173273
$ws.Range("O5").Value = 0.173378811020062
$ws.Range("P5").Value = 0.173378811020062
$ws.Range("Q5").Value = 42.65530568198744
$ws.Range("R5").Value = 383.897751137887
$ws.Range("S5").Value = 0.1031963439395094
$ws.Range("T5").Value = 0.1031963439395094

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Inhba"
$ws.Range("C6").Value = "Acvr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.433639666666666
$ws.Range("H6").Value = 25.300919
$ws.Range("I6").Value = 0.5952073574179045
$ws.Range("J6").Value = 0.5952073574179045
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.247411
$ws.Range("N6").Value = 42.742233
$ws.Range("O6").Value = 0.4883980890531961
$ws.Range("P6").Value = 0.4883980890531961
$ws.Range("Q6").Value = 120.157530556903
$ws.Range("R6").Value = 1081.417775012127
$ws.Range("S6").Value = 0.2906981359533073
$ws.Range("T6").Value = 0.2906981359533073

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Inhba"
$ws.Range("C7").Value = "Acvr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.433639666666666
$ws.Range("H7").Value = 25.300919
$ws.Range("I7").Value = 0.5952073574179045
$ws.Range("J7").Value = 0.5952073574179045
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.866548666666667
$ws.Range("N7").Value = 29.599646
$ws.Range("O7").Value = 0.3382230999267418
$ws.Range("P7").Value = 0.3382230999267418
$ws.Range("Q7").Value = 83.2109162082971
$ws.Range("R7").Value = 748.898245874674
$ws.Range("S7").Value = 0.2013128775250879
$ws.Range("T7").Value = 0.2013128775250879

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Inhba"
$ws.Range("C8").Value = "Acvr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.612402333333333
$ws.Range("H8").Value = 4.837207
$ws.Range("I8").Value = 0.1137959137276156
$ws.Range("J8").Value = 0.1137959137276156
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.057757666666666
$ws.Range("N8").Value = 15.173273
$ws.Range("O8").Value = 0.173378811020062
$ws.Range("P8").Value = 0.173378811020062
$ws.Range("Q8").Value = 8.155140263167889
$ws.Range("R8").Value = 73.396262368511
$ws.Range("S8").Value = 0.01972980022103554
$ws.Range("T8").Value = 0.01972980022103554

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Inhba"
$ws.Range("C9").Value = "Acvr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.612402333333333
$ws.Range("H9").Value = 4.837207
$ws.Range("I9").Value = 0.1137959137276156
$ws.Range("J9").Value = 0.1137959137276156
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.247411
$ws.Range("N9").Value = 42.742233
$ws.Range("O9").Value = 0.4883980890531961
$ws.Range("P9").Value = 0.4883980890531961
$ws.Range("Q9").Value = 22.972558740359
$ws.Range("R9").Value = 206.753028663231
$ws.Range("S9").Value = 0.05557770680662982
$ws.Range("T9").Value = 0.05557770680662982

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Inhba"
$ws.Range("C10").Value = "Acvr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.612402333333333
$ws.Range("H10").Value = 4.837207
$ws.Range("I10").Value = 0.1137959137276156
$ws.Range("J10").Value = 0.1137959137276156
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.866548666666667
$ws.Range("N10").Value = 29.599646
$ws.Range("O10").Value = 0.3382230999267418
$ws.Range("P10").Value = 0.3382230999267418
$ws.Range("Q10").Value = 15.90884609208022
$ws.Range("R10").Value = 143.179614828722
$ws.Range("S10").Value = 0.03848840669995022
$ws.Range("T10").Value = 0.03848840669995021
